# DataSinhVien.xlsx: "update lai khoa,lop; them nut tim kiem lop"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "anh" (photo) column J for every student row pointed at the same
# shared-string placeholder value `""`; it now holds a real (quoted) file
# path to the downloaded photo.
$ws.Range("J2:J6").Value = '"C:\Users\KHANH\Downloads\VuDinhDuc.jpg"'

# Column J ("anh") is widened to fit the new, much longer path text.
$ws.Columns.Item(10).ColumnWidth = 37

# The sheet's last saved selection moved from J6 to J8.
$ws.Range("J8").Select()
